# fix(publipostage): Try to solve Excel emoji problem
#
# The "statut" column (A) used emoji glyphs (📘 / 📕 / 📙) as status markers,
# which render unreliably in some environments (mail-merge / publipostage).
# Replace them with plain-text equivalents:
#   📘 -> ⚠️
#   📕 -> -3
#   📙 -> +3
#
# Values such as "-3"/"+3" look numeric, so a plain Value/Value2 assignment
# would make Excel store them as numbers instead of text. To keep them as
# text (matching how the original emoji were stored as shared strings) we
# assign them through a formula that evaluates to a text string, then
# collapse the formula down to its static value via copy / paste-values.
# This avoids touching the cell's number format or style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{}
$map["📘"] = "⚠️"
$map["📕"] = "-3"
$map["📙"] = "+3"

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1
$firstCol = $used.Column
$lastCol = $firstCol + $used.Columns.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        $val = $cell.Value2
        if ($null -ne $val -and $map.ContainsKey($val)) {
            $newText = $map[$val]
            $escaped = $newText.Replace('"', '""')
            $cell.Formula = '="' + $escaped + '"'
            $cell.Copy()
            $cell.PasteSpecial(-4163)  # xlPasteValues
        }
    }
}

$excel.CutCopyMode = 0
